$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the row containing account 004268684 (PATRICIA) and delete it entirely.
$patriciaRow = $ws.Range("A1").EntireColumn.Find("004268684").Row
$ws.Rows.Item($patriciaRow).Delete()

# Find the row containing account 004332783 (IRON) with the old balance 1791.89 and delete it.
$ironRow = $ws.Range("A1").EntireColumn.Find("004332783").Row
$ws.Rows.Item($ironRow).Delete()

# Find the row with GUILHERME (004224815) so we can insert the new IRON row right after it.
$guilhermeRow = $ws.Range("A1").EntireColumn.Find("004224815").Row
$insertAt = $guilhermeRow + 1

$ws.Rows.Item($insertAt).Insert()
$ws.Cells.Item($insertAt, 1).NumberFormat = "@"
$ws.Cells.Item($insertAt, 1).Value = "004332783"
$ws.Cells.Item($insertAt, 1).ClearFormats()
$ws.Cells.Item($insertAt, 2).Value = "IRON"
$ws.Cells.Item($insertAt, 3).Value = 8000
